# Generate Report for Handoff
#
# A new localization source file, dd941ab5-7374-43f6-a3c3-18528186bdff.md,
# reached "Ready for handoff" status. This inserts one new row for it on
# every sheet, directly before the existing fead0616-...md row (the new
# file sorts alphabetically right before it), pushing the
# fead0616-...md row and the trailing ".localization-config" row down by
# one. The new row mirrors the existing "Ready for handoff" row pattern.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (A:C) — File Name / zh-cn / de-de status columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(7).Insert()
$wsOverview.Range("A7").Value = "dd941ab5-7374-43f6-a3c3-18528186bdff.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"

# Rebuild hyperlinks in final left-to-right, top-to-bottom order so the
# relationship ids line up the same way Excel would renumber them.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/71428378c8763047a3507352a387fcc7af0e7406/e2e/dd941ab5-7374-43f6-a3c3-18528186bdff.md", "", "", "dd941ab5-7374-43f6-a3c3-18528186bdff.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (A:I) — per-language handoff detail.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(7).Insert()
$wsZh.Range("A7").Value = "dd941ab5-7374-43f6-a3c3-18528186bdff.md"
$wsZh.Range("B7").Value = "Ready for handoff"
$wsZh.Range("C7").Value = "dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.zh-cn.xlf"
$wsZh.Range("D7").Value = "2016-01-25 06:17:52"
$wsZh.Range("G7").Value = "0001-01-01 00:00:00"
$wsZh.Range("H7").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/288a598b1ee77e39219960e17f80572c35dcfff9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.zh-cn.xlf", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.zh-cn.xlf", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.zh-cn.xlf", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9572f2523077644a3dbe8565fbbe0ddadaeb7d8d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b474860421311042c2e35d5037bf8c00eba3310c/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61a67c1bb3d0bfc0b09208ed5006a571e79f7884/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08011f83eaae42323656547480ff840f8295f6ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.zh-cn.xlf", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/71428378c8763047a3507352a387fcc7af0e7406/e2e/dd941ab5-7374-43f6-a3c3-18528186bdff.md", "", "", "dd941ab5-7374-43f6-a3c3-18528186bdff.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49b115f46ae76283758031ddbcc8d70027d33ff8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.zh-cn.xlf", "", "", "dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d5c532270292eb43aba93eddba531e2ae101bee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.zh-cn.xlf", "", "", "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (A:I) — per-language handoff detail.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(7).Insert()
$wsDe.Range("A7").Value = "dd941ab5-7374-43f6-a3c3-18528186bdff.md"
$wsDe.Range("B7").Value = "Ready for handoff"
$wsDe.Range("C7").Value = "dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.de-de.xlf"
$wsDe.Range("D7").Value = "2016-01-25 06:18:03"
$wsDe.Range("G7").Value = "0001-01-01 00:00:00"
$wsDe.Range("H7").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91e1e715d25ad71db3ba3ce7d1667df85facedfa/e2e/4f17c1f2-6627-40b3-b43d-91c7169672d9.md", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/288a598b1ee77e39219960e17f80572c35dcfff9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.de-de.xlf", "", "", "4f17c1f2-6627-40b3-b43d-91c7169672d9.c955d1b25919dc449426e0a90610e6375ad79b2d.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/56e46f49-4e13-4895-8960-5bb9e3990598.md", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.de-de.xlf", "", "", "56e46f49-4e13-4895-8960-5bb9e3990598.972d26d186fafdcb8dde947c8ee7ae69178b62dd.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/220c142ea87d40a737516fea595ac1faff88d1df/e2e/cbc66d11-eedb-4924-9ea8-e10b3ffda301.md", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b687909adcb62f1c3da5a2a47a11d2bea3b80f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.de-de.xlf", "", "", "cbc66d11-eedb-4924-9ea8-e10b3ffda301.066968cb1bedba9e66f2ebad52f0b867a9348951.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2e2a93a0809622b70486736ced2d76133905d8/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9572f2523077644a3dbe8565fbbe0ddadaeb7d8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b474860421311042c2e35d5037bf8c00eba3310c/e2e/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/61a67c1bb3d0bfc0b09208ed5006a571e79f7884/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf", "", "", "e61a1424-bd49-45bb-a42e-1dc45ef1bb80.0e3571501821ba00efd607907b63b28138246024.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/14974f658b0d3a1cf28a3a58919cf22aeedc75b1/e2e/268b50a0-f412-4f69-99e2-079bdfdf1585.md", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08011f83eaae42323656547480ff840f8295f6ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.de-de.xlf", "", "", "268b50a0-f412-4f69-99e2-079bdfdf1585.82626d270fb5faec983355aec6122bf59d5b1010.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/71428378c8763047a3507352a387fcc7af0e7406/e2e/dd941ab5-7374-43f6-a3c3-18528186bdff.md", "", "", "dd941ab5-7374-43f6-a3c3-18528186bdff.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ae4fe4921e078f4089969457ba61e4c334b5e10/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.de-de.xlf", "", "", "dd941ab5-7374-43f6-a3c3-18528186bdff.cd2bc6fe53a621ab70db369f3363e96d2edd8c52.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/e2e/fead0616-771c-481d-b56a-5918df0efd59.md", "", "", "fead0616-771c-481d-b56a-5918df0efd59.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03bca571b259ae7863be9e4345aff163dc6b69eb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.de-de.xlf", "", "", "fead0616-771c-481d-b56a-5918df0efd59.74e64ac24ffce7269733ec7ed43ca6329d62d70a.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/16ff55d06f0a3c520c77e066d353a9cfa71f8579/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Added dd941ab5-7374-43f6-a3c3-18528186bdff.md as Ready for handoff on all sheets."
